# Remove the footer block that used to follow the "LOQ4205: Sistemas
# Produtivos II (Requisito fraco)" requirement line:
#   - the blank paragraph right after it
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#     pages. Original theme under Creative Commons Attribution"

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter..." paragraph via Find.
$findRange = $d.Content
[void]$findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$startParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Start -le $findRange.Start -and $findRange.Start -lt $para.Range.End) {
        $startParaIndex = $i
        break
    }
}

# The empty paragraph right before "Ver no Jupiter..." and the copyright
# paragraph right after it bound the block to remove.
$blankPara = $d.Paragraphs.Item($startParaIndex - 1)
$copyrightPara = $d.Paragraphs.Item($startParaIndex + 1)

$deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
$deleteRange.Delete()
